# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) on several per-job leve sheets to reflect the latest pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 376327.56
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 376327.56
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1128982.68
$ws.Range("N17").Value = -1129318.68

$ws.Range("H41").Value = 760.4286
$ws.Range("I41").Value = 787
$ws.Range("J41").Value = 725
$ws.Range("K41").Value = 787
$ws.Range("L41").Value = 725
$ws.Range("M41").Value = -347
$ws.Range("N41").Value = -1605

$ws.Range("H64").Value = 4353.3887
$ws.Range("I64").Value = 3880
$ws.Range("J64").Value = 4535.4614
$ws.Range("K64").Value = 3880
$ws.Range("L64").Value = 4535.4614
$ws.Range("M64").Value = -3632
$ws.Range("N64").Value = -5031.4614

$ws.Range("H67").Value = 4353.3887
$ws.Range("I67").Value = 3880
$ws.Range("J67").Value = 4535.4614
$ws.Range("K67").Value = 3880
$ws.Range("L67").Value = 4535.4614
$ws.Range("M67").Value = -3022
$ws.Range("N67").Value = -6251.4614

$ws.Range("H86").Value = 2207.4666
$ws.Range("I86").Value = 1866.6666
$ws.Range("J86").Value = 2292.6667
$ws.Range("K86").Value = 1866.6666
$ws.Range("L86").Value = 2292.6667
$ws.Range("M86").Value = -743.6666
$ws.Range("N86").Value = -4538.6667

$ws.Range("H89").Value = 2207.4666
$ws.Range("I89").Value = 1866.6666
$ws.Range("J89").Value = 2292.6667
$ws.Range("K89").Value = 9333.333000000001
$ws.Range("L89").Value = 11463.3335
$ws.Range("M89").Value = -3717.333000000001
$ws.Range("N89").Value = -22695.3335

$ws.Range("H107").Value = 1168.6666
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 1653
$ws.Range("K107").Value = 200
$ws.Range("L107").Value = 1653
$ws.Range("M107").Value = 1720

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6600.7646
$ws.Range("I88").Value = 2996.5
$ws.Range("J88").Value = 8566.727999999999
$ws.Range("K88").Value = 2996.5
$ws.Range("L88").Value = 8566.727999999999
$ws.Range("M88").Value = -2590.5
$ws.Range("N88").Value = -9378.727999999999

$ws.Range("H91").Value = 6600.7646
$ws.Range("I91").Value = 2996.5
$ws.Range("J91").Value = 8566.727999999999
$ws.Range("K91").Value = 2996.5
$ws.Range("L91").Value = 8566.727999999999
$ws.Range("M91").Value = -1592.5
$ws.Range("N91").Value = -11374.728

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()

$ws.Range("H124").Value = 28189.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 28189.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 28189.5
$ws.Range("N124").Value = -38009.5

$ws.Range("H125").Value = 60000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 60000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

$ws.Range("H127").Value = 30000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 30000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920

$ws.Range("H128").Value = 44000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 44000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 44000
$ws.Range("N128").Value = -53960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18343.572
$ws.Range("I86").Value = 26932.889
$ws.Range("J86").Value = 2882.8
$ws.Range("K86").Value = 26932.889
$ws.Range("L86").Value = 2882.8
$ws.Range("M86").Value = -25809.889
$ws.Range("N86").Value = -5128.8

$ws.Range("H89").Value = 18343.572
$ws.Range("I89").Value = 26932.889
$ws.Range("J89").Value = 2882.8
$ws.Range("K89").Value = 134664.445
$ws.Range("L89").Value = 14414
$ws.Range("M89").Value = -129048.445
$ws.Range("N89").Value = -25646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 135.42857
$ws.Range("I19").Value = 135.42857
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 135.42857
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 34.57142999999999

$ws.Range("H24").Value = 135.42857
$ws.Range("I24").Value = 135.42857
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 135.42857
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 34.57142999999999

$ws.Range("H62").Value = 3809.625
$ws.Range("I62").Value = 3300.8333
$ws.Range("J62").Value = 5336
$ws.Range("K62").Value = 3300.8333
$ws.Range("L62").Value = 5336
$ws.Range("M62").Value = -2676.8333

$ws.Range("H65").Value = 3809.625
$ws.Range("I65").Value = 3300.8333
$ws.Range("J65").Value = 5336
$ws.Range("K65").Value = 16504.1665
$ws.Range("L65").Value = 26680
$ws.Range("M65").Value = -13384.1665

$ws.Range("H95").Value = 20000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 20000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.8
$ws.Range("I2").Value = 69.59999999999999
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 417.6
$ws.Range("L2").Value = 156
$ws.Range("M2").Value = -304.6
$ws.Range("N2").Value = -382

$ws.Range("H23").Value = 162.8
$ws.Range("I23").Value = 170
$ws.Range("J23").Value = 161
$ws.Range("K23").Value = 510
$ws.Range("L23").Value = 483
$ws.Range("M23").Value = -275
$ws.Range("N23").Value = -953

$ws.Range("H38").Value = 148.08
$ws.Range("I38").Value = 45.857143
$ws.Range("J38").Value = 187.83333
$ws.Range("K38").Value = 137.571429
$ws.Range("L38").Value = 563.49999
$ws.Range("M38").Value = 209.428571
$ws.Range("N38").Value = -1257.49999

$ws.Range("H132").Value = 1275.0526
$ws.Range("I132").Value = 940.46155
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 8464.15395
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -5934.15395
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1668.1666
$ws.Range("I102").Value = 1531.9615
$ws.Range("J102").Value = 2553.5
$ws.Range("K102").Value = 1531.9615
$ws.Range("L102").Value = 2553.5
$ws.Range("M102").Value = 90.03850000000011
$ws.Range("N102").Value = -5797.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 2458.2
$ws.Range("I81").Value = 837.5
$ws.Range("J81").Value = 2863.375
$ws.Range("K81").Value = 1675
$ws.Range("L81").Value = 5726.75
$ws.Range("M81").Value = -614
$ws.Range("N81").Value = -7848.75

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 2458.2
$ws.Range("I84").Value = 837.5
$ws.Range("J84").Value = 2863.375
$ws.Range("K84").Value = 8375
$ws.Range("L84").Value = 28633.75
$ws.Range("M84").Value = -3071
$ws.Range("N84").Value = -39241.75

